$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 108.3097613333333
$ws.Range("H2").Value = 324.929284
$ws.Range("I2").Value = 0.9760647858278649
$ws.Range("J2").Value = 0.9760647858278649
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 12677.0686838943
$ws.Range("R2").Value = 114093.6181550487
$ws.Range("S2").Value = 0.316768712363761
$ws.Range("T2").Value = 0.316768712363761
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 108.3097613333333
$ws.Range("H3").Value = 324.929284
$ws.Range("I3").Value = 0.9760647858278649
$ws.Range("J3").Value = 0.9760647858278649
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 11002.10959980442
$ws.Range("R3").Value = 99018.98639823981
$ws.Range("S3").Value = 0.2749156116541934
$ws.Range("T3").Value = 0.2749156116541934
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 108.3097613333333
$ws.Range("H4").Value = 324.929284
$ws.Range("I4").Value = 0.9760647858278649
$ws.Range("J4").Value = 0.9760647858278649
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 15382.88765563295
$ws.Range("R4").Value = 138445.9889006965
$ws.Range("S4").Value = 0.3843804618099104
$ws.Range("T4").Value = 0.3843804618099104
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.160250666666667
$ws.Range("H5").Value = 6.480752000000001
$ws.Range("I5").Value = 0.0194677245922947
$ws.Range("J5").Value = 0.0194677245922947
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 252.8455952504587
$ws.Range("R5").Value = 2275.610357254128
$ws.Range("S5").Value = 0.006317988458648342
$ws.Range("T5").Value = 0.006317988458648342
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.160250666666667
$ws.Range("H6").Value = 6.480752000000001
$ws.Range("I6").Value = 0.0194677245922947
$ws.Range("J6").Value = 0.0194677245922947
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 219.4383433693583
$ws.Range("R6").Value = 1974.945090324224
$ws.Range("S6").Value = 0.0054832235436777
$ws.Range("T6").Value = 0.005483223543677699
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.160250666666667
$ws.Range("H7").Value = 6.480752000000001
$ws.Range("I7").Value = 0.0194677245922947
$ws.Range("J7").Value = 0.0194677245922947
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 306.8134663418596
$ws.Range("R7").Value = 2761.321197076737
$ws.Range("S7").Value = 0.007666512589968656
$ws.Range("T7").Value = 0.007666512589968656
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4957383333333333
$ws.Range("H8").Value = 1.487215
$ws.Range("I8").Value = 0.004467489579840358
$ws.Range("J8").Value = 0.004467489579840358
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 58.02347658734833
$ws.Range("R8").Value = 522.211289286135
$ws.Range("S8").Value = 0.001449863720372064
$ws.Range("T8").Value = 0.001449863720372064
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4957383333333333
$ws.Range("H9").Value = 1.487215
$ws.Range("I9").Value = 0.004467489579840358
$ws.Range("J9").Value = 0.004467489579840358
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 50.35711840756445
$ws.Range("R9").Value = 453.21406566808
$ws.Range("S9").Value = 0.00125830031800486
$ws.Range("T9").Value = 0.001258300318004859
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4957383333333333
$ws.Range("H10").Value = 1.487215
$ws.Range("I10").Value = 0.004467489579840358
$ws.Range("J10").Value = 0.004467489579840358
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 70.40812383279112
$ws.Range("R10").Value = 633.67311449512
$ws.Range("S10").Value = 0.001759325541463434
$ws.Range("T10").Value = 0.001759325541463434

Write-Output "Applied Natmi following Dr Hou advice edits"
